$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet11")

# Header row (row 6): rotate values left across B,C,D (B<-C, C<-D, D<-B)
$ws.Range("B6").Value = "does"
$ws.Range("C6").Value = "this"
$ws.Range("D6").Value = "work"

# Row 7
$ws.Range("B7").Value = "[0, 0]"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2

# Row 8
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 2

# Row 9
$ws.Range("B9").Value = "[0, 0]"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2

# Row 10
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
